$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 222-224; Excel shifts existing rows 222:265 down to 225:268
# and copies the formatting (e.g. the date number format on column D) from the row above.
$ws.Rows("222:224").Insert()

# New row 222 - Calidad "Especial", Fecha 2021-11-04 (serial 44504)
$ws.Cells.Item(222, 1).Value = 2
$ws.Cells.Item(222, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(222, 3).Value = "Coquimbo"
$ws.Cells.Item(222, 4).Value = 44504
$ws.Cells.Item(222, 5).Value = 4
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100101
$ws.Cells.Item(222, 8).Value = "Berries"
$ws.Cells.Item(222, 9).Value = 100112025
$ws.Cells.Item(222, 10).Value = "Frutilla"
$ws.Cells.Item(222, 11).Value = "Sin especificar"
$ws.Cells.Item(222, 12).Value = "Especial"
$ws.Cells.Item(222, 13).Value = 400
$ws.Cells.Item(222, 14).Value = 12500
$ws.Cells.Item(222, 15).Value = 13000
$ws.Cells.Item(222, 16).Value = 12750
$ws.Cells.Item(222, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(222, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(222, 19).Value = 1821
$ws.Cells.Item(222, 20).Value = 7

# New row 223 - Calidad "Primera", Fecha 2021-11-04 (serial 44504)
$ws.Cells.Item(223, 1).Value = 2
$ws.Cells.Item(223, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(223, 3).Value = "Coquimbo"
$ws.Cells.Item(223, 4).Value = 44504
$ws.Cells.Item(223, 5).Value = 4
$ws.Cells.Item(223, 6).Value = "Fruta"
$ws.Cells.Item(223, 7).Value = 100101
$ws.Cells.Item(223, 8).Value = "Berries"
$ws.Cells.Item(223, 9).Value = 100112025
$ws.Cells.Item(223, 10).Value = "Frutilla"
$ws.Cells.Item(223, 11).Value = "Sin especificar"
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 340
$ws.Cells.Item(223, 14).Value = 10500
$ws.Cells.Item(223, 15).Value = 11000
$ws.Cells.Item(223, 16).Value = 10750
$ws.Cells.Item(223, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(223, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(223, 19).Value = 1536
$ws.Cells.Item(223, 20).Value = 7

# New row 224 - Calidad "Segunda", Fecha 2021-11-04 (serial 44504)
$ws.Cells.Item(224, 1).Value = 2
$ws.Cells.Item(224, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44504
$ws.Cells.Item(224, 5).Value = 4
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100112025
$ws.Cells.Item(224, 10).Value = "Frutilla"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Segunda"
$ws.Cells.Item(224, 13).Value = 240
$ws.Cells.Item(224, 14).Value = 8500
$ws.Cells.Item(224, 15).Value = 9000
$ws.Cells.Item(224, 16).Value = 8750
$ws.Cells.Item(224, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(224, 19).Value = 1250
$ws.Cells.Item(224, 20).Value = 7
